$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "工作表2"
Write-Output $wb.Worksheets.Count
Write-Output $wb.Worksheets.Item(1).Name
Write-Output $wb.Worksheets.Item(2).Name
